$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Step 1: delete the three rows that held 0.00018 / 0.00021 / 0.00025
#     (originally rows 9, 10, 11) ---
$t.Rows.Item(9).Delete()
$t.Rows.Item(9).Delete()
$t.Rows.Item(9).Delete()

# --- Step 2: insert three new rows right before the (now-shifted) row 5
#     (the row holding 0.00004). Insert in reverse so the final order
#     reads 0.00003, 0.00059, 0.00019. ---
$newRowC = $t.Rows.Add($t.Rows.Item(5))
$newRowC.Cells.Item(1).Range.Text = "0.00019"

$newRowB = $t.Rows.Add($t.Rows.Item(5))
$newRowB.Cells.Item(1).Range.Text = "0.00059"

$newRowA = $t.Rows.Add($t.Rows.Item(5))
$newRowA.Cells.Item(1).Range.Text = "0.00003"

# --- Step 3: simple text replacements on the fixed leading rows ---
$t.Cell(1,1).Range.Text = "0M"
$t.Cell(2,1).Range.Text = "0M"
$t.Cell(3,1).Range.Text = "0M"
$t.Cell(4,1).Range.Text = "360"

# --- Step 4: text replacements on the rows following the inserted block
#     (originally rows 6, 7, 8, 12 -> now rows 9, 10, 11, 12) ---
$t.Cell(9,1).Range.Text = "0.00040"
$t.Cell(10,1).Range.Text = "0.00044"
$t.Cell(11,1).Range.Text = "0.00047"
$t.Cell(12,1).Range.Text = "0.08104"

# --- Step 5: collapse the three multi-run summary rows at the bottom
#     (rows 44, 45, 46) down to their single leading value ---
$t.Cell(44,1).Range.Text = "99.96"
$t.Cell(45,1).Range.Text = "0.08"
$t.Cell(46,1).Range.Text = "210"
